$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6132958801498127
$ws1.Range("C2").Value = 0.5645677694770544
$ws1.Range("D2").Value = 0.9906367041198502
$ws1.Range("E2").Value = 0.7192386131883073
$ws1.Range("F2").Value = 0.8607224210868858
$ws1.Range("G2").Value = 0.9626933575978162
$ws1.Range("H2").Value = 0.7775954214535201
$ws1.Range("I2").Value = 529
$ws1.Range("J2").Value = 408
$ws1.Range("K2").Value = 126
$ws1.Range("L2").Value = 5

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 - class "0"
$ws2.Range("B2").Value = 0.9618320610687023
$ws2.Range("C2").Value = 0.2359550561797753
$ws2.Range("D2").Value = 0.3789473684210526

# Row 3 - class "1"
$ws2.Range("B3").Value = 0.5645677694770544
$ws2.Range("C3").Value = 0.9906367041198502
$ws2.Range("D3").Value = 0.7192386131883073

# Row 4 - accuracy
$ws2.Range("B4").Value = 0.6132958801498127
$ws2.Range("C4").Value = 0.6132958801498127
$ws2.Range("D4").Value = 0.6132958801498127
$ws2.Range("E4").Value = 0.6132958801498127

# Row 5 - macro avg
$ws2.Range("B5").Value = 0.7631999152728783
$ws2.Range("C5").Value = 0.6132958801498127
$ws2.Range("D5").Value = 0.54909299080468

# Row 6 - weighted avg
$ws2.Range("B6").Value = 0.7631999152728784
$ws2.Range("C6").Value = 0.6132958801498127
$ws2.Range("D6").Value = 0.54909299080468

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - Actual 0
$ws3.Range("B2").Value = 126
$ws3.Range("C2").Value = 408

# Row 3 - Actual 1
$ws3.Range("B3").Value = 5
$ws3.Range("C3").Value = 529
